$wb = $excel.ActiveWorkbook

# Sheet 1: departements
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(681, 3).Value = 5.217391304347826
$ws1.Cells.Item(681, 4).Value = 6
$ws1.Cells.Item(682, 3).Value = 5.555555555555555
$ws1.Cells.Item(682, 4).Value = 5
$ws1.Cells.Item(683, 3).Value = 2.531645569620253
$ws1.Cells.Item(683, 4).Value = 2
$ws1.Cells.Item(683, 5).Value = 79
$ws1.Cells.Item(686, 3).Value = 3.424657534246575
$ws1.Cells.Item(686, 4).Value = 5
$ws1.Cells.Item(687, 3).Value = 1.449275362318841
$ws1.Cells.Item(687, 4).Value = 1
$ws1.Cells.Item(688, 3).Value = 4.081632653061225
$ws1.Cells.Item(688, 4).Value = 2
$ws1.Cells.Item(690, 3).Value = 1.538461538461539
$ws1.Cells.Item(690, 4).Value = 1
$ws1.Cells.Item(693, 3).Value = 4.014598540145985
$ws1.Cells.Item(693, 4).Value = 11
$ws1.Cells.Item(693, 5).Value = 274
$ws1.Cells.Item(694, 3).Value = 1.666666666666667
$ws1.Cells.Item(694, 4).Value = 2
$ws1.Cells.Item(695, 3).Value = 1.851851851851852
$ws1.Cells.Item(695, 5).Value = 54
$ws1.Cells.Item(698, 3).Value = 4.918032786885246
$ws1.Cells.Item(698, 4).Value = 3
$ws1.Cells.Item(700, 3).Value = 1.626016260162602
$ws1.Cells.Item(700, 4).Value = 2
$ws1.Cells.Item(701, 3).Value = 2.5
$ws1.Cells.Item(701, 4).Value = 3
$ws1.Cells.Item(702, 3).Value = 3.333333333333333
$ws1.Cells.Item(702, 5).Value = 30
$ws1.Cells.Item(703, 3).Value = 1.08695652173913
$ws1.Cells.Item(703, 4).Value = 1
$ws1.Cells.Item(704, 3).Value = 2
$ws1.Cells.Item(704, 4).Value = 2
$ws1.Cells.Item(704, 5).Value = 100
$ws1.Cells.Item(706, 3).Value = 4.210526315789473
$ws1.Cells.Item(706, 4).Value = 4
$ws1.Cells.Item(708, 3).Value = 1.36986301369863
$ws1.Cells.Item(708, 4).Value = 2
$ws1.Cells.Item(711, 3).Value = 2.290076335877862
$ws1.Cells.Item(711, 4).Value = 3
$ws1.Cells.Item(711, 5).Value = 131
$ws1.Cells.Item(712, 3).Value = 6.030150753768844
$ws1.Cells.Item(712, 4).Value = 12
$ws1.Cells.Item(715, 3).Value = 3.164556962025316
$ws1.Cells.Item(715, 4).Value = 5
$ws1.Cells.Item(716, 3).Value = 1.351351351351351
$ws1.Cells.Item(716, 4).Value = 2
$ws1.Cells.Item(717, 3).Value = 4.081632653061225
$ws1.Cells.Item(717, 4).Value = 2
$ws1.Cells.Item(717, 5).Value = 49
$ws1.Cells.Item(718, 3).Value = 5.263157894736842
$ws1.Cells.Item(718, 4).Value = 5
$ws1.Cells.Item(719, 3).Value = 3.783783783783784
$ws1.Cells.Item(719, 4).Value = 7
$ws1.Cells.Item(719, 5).Value = 185
$ws1.Cells.Item(720, 3).Value = 1.298701298701299
$ws1.Cells.Item(720, 5).Value = 77
$ws1.Cells.Item(722, 3).Value = 6.349206349206349
$ws1.Cells.Item(722, 4).Value = 4
$ws1.Cells.Item(723, 3).Value = 4.615384615384616
$ws1.Cells.Item(723, 4).Value = 6
$ws1.Cells.Item(723, 5).Value = 130
$ws1.Cells.Item(724, 5).Value = 52
$ws1.Cells.Item(725, 3).Value = 1.704545454545454
$ws1.Cells.Item(725, 5).Value = 176
$ws1.Cells.Item(726, 3).Value = 1.834862385321101
$ws1.Cells.Item(726, 5).Value = 109
$ws1.Cells.Item(727, 3).Value = 3.389830508474576
$ws1.Cells.Item(727, 4).Value = 2
$ws1.Cells.Item(729, 5).Value = 37
$ws1.Cells.Item(730, 3).Value = 1.851851851851852
$ws1.Cells.Item(730, 4).Value = 2
$ws1.Cells.Item(731, 3).Value = 0.847457627118644
$ws1.Cells.Item(731, 4).Value = 1
$ws1.Cells.Item(732, 3).Value = 6.382978723404255
$ws1.Cells.Item(732, 4).Value = 6
$ws1.Cells.Item(732, 5).Value = 94
$ws1.Cells.Item(733, 3).Value = 1.96078431372549
$ws1.Cells.Item(733, 5).Value = 51
$ws1.Cells.Item(735, 3).Value = 3.225806451612903
$ws1.Cells.Item(735, 4).Value = 3
$ws1.Cells.Item(736, 3).Value = 2.857142857142857
$ws1.Cells.Item(736, 4).Value = 1
$ws1.Cells.Item(737, 3).Value = 3.571428571428571
$ws1.Cells.Item(737, 4).Value = 5
$ws1.Cells.Item(740, 3).Value = 5.2
$ws1.Cells.Item(740, 4).Value = 13
$ws1.Cells.Item(740, 5).Value = 250
$ws1.Cells.Item(741, 3).Value = 7.627118644067797
$ws1.Cells.Item(741, 4).Value = 9
$ws1.Cells.Item(743, 3).Value = 7.327586206896551
$ws1.Cells.Item(743, 4).Value = 17
$ws1.Cells.Item(743, 5).Value = 232
$ws1.Cells.Item(744, 3).Value = 0.8403361344537815
$ws1.Cells.Item(744, 4).Value = 1
$ws1.Cells.Item(747, 3).Value = 1.449275362318841
$ws1.Cells.Item(747, 4).Value = 1
$ws1.Cells.Item(748, 3).Value = 5.405405405405405
$ws1.Cells.Item(748, 5).Value = 148
$ws1.Cells.Item(749, 3).Value = 1.739130434782609
$ws1.Cells.Item(749, 4).Value = 2
$ws1.Cells.Item(750, 3).Value = 5.5
$ws1.Cells.Item(750, 4).Value = 11
$ws1.Cells.Item(751, 3).Value = 11.86440677966102
$ws1.Cells.Item(751, 5).Value = 59
$ws1.Cells.Item(752, 3).Value = 2.836879432624114
$ws1.Cells.Item(752, 5).Value = 141
$ws1.Cells.Item(753, 3).Value = 3.96039603960396
$ws1.Cells.Item(753, 4).Value = 4
$ws1.Cells.Item(756, 3).Value = 9.803921568627452
$ws1.Cells.Item(756, 4).Value = 5
$ws1.Cells.Item(758, 3).Value = 8.426966292134832
$ws1.Cells.Item(758, 4).Value = 15
$ws1.Cells.Item(758, 5).Value = 178
$ws1.Cells.Item(759, 3).Value = 8.783783783783784
$ws1.Cells.Item(759, 4).Value = 13
$ws1.Cells.Item(761, 3).Value = 6.25
$ws1.Cells.Item(761, 4).Value = 7
$ws1.Cells.Item(761, 5).Value = 112
$ws1.Cells.Item(762, 5).Value = 71
$ws1.Cells.Item(763, 3).Value = 3.125
$ws1.Cells.Item(763, 4).Value = 2
$ws1.Cells.Item(764, 3).Value = 5.405405405405405
$ws1.Cells.Item(764, 4).Value = 10
$ws1.Cells.Item(766, 3).Value = 0.8547008547008548
$ws1.Cells.Item(766, 4).Value = 1
$ws1.Cells.Item(767, 3).Value = 1.351351351351351
$ws1.Cells.Item(767, 4).Value = 1
$ws1.Cells.Item(768, 3).Value = 1.25
$ws1.Cells.Item(768, 4).Value = 1
$ws1.Cells.Item(770, 3).Value = 5.747126436781609
$ws1.Cells.Item(770, 4).Value = 5
$ws1.Cells.Item(772, 3).Value = 13.76811594202899
$ws1.Cells.Item(772, 4).Value = 19
$ws1.Cells.Item(774, 3).Value = 16.48351648351648
$ws1.Cells.Item(774, 4).Value = 15
$ws1.Cells.Item(776, 3).Value = 14.15929203539823
$ws1.Cells.Item(776, 4).Value = 16

# Sheet 2: regions
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(8, 4).Value = 10.97
$ws2.Cells.Item(8, 5).Value = 97
$ws2.Cells.Item(8, 6).Value = 884
$ws2.Cells.Item(17, 4).Value = 3.78
$ws2.Cells.Item(17, 5).Value = 17
$ws2.Cells.Item(17, 6).Value = 450
$ws2.Cells.Item(26, 4).Value = 3.31
$ws2.Cells.Item(26, 5).Value = 22
$ws2.Cells.Item(26, 6).Value = 664
$ws2.Cells.Item(35, 4).Value = 1.96
$ws2.Cells.Item(35, 5).Value = 11
$ws2.Cells.Item(44, 4).Value = 6.36
$ws2.Cells.Item(44, 5).Value = 51
$ws2.Cells.Item(44, 6).Value = 802
$ws2.Cells.Item(53, 4).Value = 3.36
$ws2.Cells.Item(53, 5).Value = 29
$ws2.Cells.Item(53, 6).Value = 863
$ws2.Cells.Item(62, 4).Value = 2.15
$ws2.Cells.Item(62, 5).Value = 12
$ws2.Cells.Item(62, 6).Value = 558
$ws2.Cells.Item(71, 4).Value = 2.17
$ws2.Cells.Item(71, 5).Value = 12
$ws2.Cells.Item(80, 4).Value = 1.52
$ws2.Cells.Item(80, 5).Value = 17
$ws2.Cells.Item(80, 6).Value = 1118
$ws2.Cells.Item(89, 4).Value = 2.83
$ws2.Cells.Item(89, 5).Value = 31
$ws2.Cells.Item(89, 6).Value = 1094
$ws2.Cells.Item(98, 4).Value = 3.08
$ws2.Cells.Item(98, 5).Value = 41
$ws2.Cells.Item(98, 6).Value = 1331
$ws2.Cells.Item(107, 4).Value = 3.47
$ws2.Cells.Item(107, 5).Value = 28
$ws2.Cells.Item(107, 6).Value = 806

# Sheet 3: national
$ws3 = $wb.Worksheets.Item(3)
$ws3.Cells.Item(8, 2).Value = 3.75
$ws3.Cells.Item(8, 3).Value = 368
$ws3.Cells.Item(8, 4).Value = 9816
